$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B18: change from text "4" to numeric 4
$ws.Cells.Item(18, 2).Value = 4

# Add new row 19
$ws.Cells.Item(19, 1).Value = "Ruilin"
$ws.Cells.Item(19, 2).Value = "'3"
$ws.Cells.Item(19, 2).Style = "Normal"
$ws.Cells.Item(19, 3).Value = "无"
$ws.Cells.Item(19, 4).Value = "DIS"
$ws.Cells.Item(19, 5).Value = "OTH"
$ws.Cells.Item(19, 6).Value = "2bb8b329-99fa-4c06-a5b4-7897e3cce401"
$ws.Cells.Item(19, 7).Value = "S1PWi_lC-_annotated.xlsx"
$ws.Cells.Item(19, 8).Value = "Each network is trained with 50 epochs."
